$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 941
$ws.Range("I98").Value = 1005.2
$ws.Range("J98").Value = 620
$ws.Range("K98").Value = 1005.2
$ws.Range("L98").Value = 620
$ws.Range("M98").Value = 492.8
$ws.Range("N98").Value = -3616

$ws.Range("H122").Value = 941
$ws.Range("I122").Value = 1005.2
$ws.Range("J122").Value = 620
$ws.Range("K122").Value = 3015.6
$ws.Range("L122").Value = 1860
$ws.Range("M122").Value = -565.6000000000004
$ws.Range("N122").Value = -6760

$ws.Range("H137").Value = 1177.75
$ws.Range("I137").Value = 1177.75
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3533.25
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -983.25
$ws.Range("N137").ClearContents()


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1874006.6
$ws.Range("I32").Value = 3833.0344
$ws.Range("J32").Value = 13926236
$ws.Range("K32").Value = 3833.0344
$ws.Range("L32").Value = 13926236
$ws.Range("M32").Value = -3546.0344
$ws.Range("N32").Value = -13926810

$ws.Range("H45").Value = 2559.6
$ws.Range("I45").Value = 1949.5
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 1949.5
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1572.5
$ws.Range("N45").Value = -5754

$ws.Range("H61").Value = 1408.1628
$ws.Range("I61").Value = 1263.4688
$ws.Range("K61").Value = 1263.4688
$ws.Range("M61").Value = -1051.4688

$ws.Range("H74").Value = 795.6491
$ws.Range("I74").Value = 783.0682
$ws.Range("J74").Value = 838.2308
$ws.Range("K74").Value = 783.0682
$ws.Range("L74").Value = 838.2308
$ws.Range("M74").Value = 90.93179999999995
$ws.Range("N74").Value = -2586.2308

$ws.Range("H77").Value = 795.6491
$ws.Range("I77").Value = 783.0682
$ws.Range("J77").Value = 838.2308
$ws.Range("K77").Value = 3915.341
$ws.Range("L77").Value = 4191.154
$ws.Range("M77").Value = 452.6589999999997
$ws.Range("N77").Value = -12927.154

$ws.Range("H88").Value = 3350
$ws.Range("I88").Value = 2800
$ws.Range("K88").Value = 2800
$ws.Range("M88").Value = -2394

$ws.Range("H91").Value = 3350
$ws.Range("I91").Value = 2800
$ws.Range("K91").Value = 2800
$ws.Range("M91").Value = -1396

$ws.Range("H97").Value = 541.12
$ws.Range("I97").Value = 517.3158
$ws.Range("K97").Value = 517.3158
$ws.Range("M97").Value = -21.31579999999997

$ws.Range("H110").Value = 404.5
$ws.Range("I110").Value = 434.57144
$ws.Range("K110").Value = 434.57144
$ws.Range("M110").Value = 1610.42856

$ws.Range("H136").Value = 1408.1628
$ws.Range("I136").Value = 1263.4688
$ws.Range("K136").Value = 3790.4064
$ws.Range("M136").Value = -1240.4064


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2894.4443
$ws.Range("I86").Value = 2455.5557
$ws.Range("J86").Value = 3333.3333
$ws.Range("K86").Value = 2455.5557
$ws.Range("L86").Value = 3333.3333
$ws.Range("M86").Value = -1332.5557
$ws.Range("N86").Value = -5579.3333

$ws.Range("H89").Value = 2894.4443
$ws.Range("I89").Value = 2455.5557
$ws.Range("J89").Value = 3333.3333
$ws.Range("K89").Value = 12277.7785
$ws.Range("L89").Value = 16666.6665
$ws.Range("M89").Value = -6661.7785
$ws.Range("N89").Value = -27898.6665

$ws.Range("H107").Value = 19705.857
$ws.Range("I107").Value = 856.8333
$ws.Range("J107").Value = 132800
$ws.Range("K107").Value = 856.8333
$ws.Range("L107").Value = 132800
$ws.Range("M107").Value = 1063.1667
$ws.Range("N107").Value = -136640

$ws.Range("H134").Value = 56934.38
$ws.Range("I134").Value = 3172.3447
$ws.Range("J134").Value = 251821.75
$ws.Range("K134").Value = 9517.034100000001
$ws.Range("L134").Value = 755465.25
$ws.Range("M134").Value = -6982.034100000001
$ws.Range("N134").Value = -760535.25


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2768.913
$ws.Range("I31").Value = 2768.913
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2768.913
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2473.913
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 2768.913
$ws.Range("I34").Value = 2768.913
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2768.913
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2566.913
$ws.Range("N34").ClearContents()

$ws.Range("H58").Value = 7242.5264
$ws.Range("I58").Value = 2769
$ws.Range("J58").Value = 10496
$ws.Range("K58").Value = 2769
$ws.Range("L58").Value = 10496
$ws.Range("M58").Value = -2566
$ws.Range("N58").Value = -10902

$ws.Range("H62").Value = 10171.429
$ws.Range("I62").Value = 11550
$ws.Range("J62").Value = 8333.333000000001
$ws.Range("K62").Value = 11550
$ws.Range("L62").Value = 8333.333000000001
$ws.Range("M62").Value = -10926
$ws.Range("N62").Value = -9581.333000000001

$ws.Range("H65").Value = 10171.429
$ws.Range("I65").Value = 11550
$ws.Range("J65").Value = 8333.333000000001
$ws.Range("K65").Value = 57750
$ws.Range("L65").Value = 41666.665
$ws.Range("M65").Value = -54630
$ws.Range("N65").Value = -47906.665

$ws.Range("H134").Value = 1813.8235
$ws.Range("I134").Value = 1604.9584
$ws.Range("J134").Value = 2315.1
$ws.Range("K134").Value = 4814.8752
$ws.Range("L134").Value = 6945.299999999999
$ws.Range("M134").Value = -2279.8752
$ws.Range("N134").Value = -12015.3

$ws.Range("H136").Value = 7242.5264
$ws.Range("I136").Value = 2769
$ws.Range("J136").Value = 10496
$ws.Range("K136").Value = 8307
$ws.Range("L136").Value = 31488
$ws.Range("M136").Value = -5757
$ws.Range("N136").Value = -36588


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5113669
$ws.Range("I131").Value = 23857958
$ws.Range("J131").Value = 1590.5714
$ws.Range("K131").Value = 71573874
$ws.Range("L131").Value = 4771.7142
$ws.Range("M131").Value = -71568834
$ws.Range("N131").Value = -14851.7142


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4626.125
$ws.Range("I70").Value = 4161.5
$ws.Range("J70").Value = 6020
$ws.Range("K70").Value = 4161.5
$ws.Range("L70").Value = 6020
$ws.Range("M70").Value = -3891.5
$ws.Range("N70").Value = -6560

$ws.Range("H73").Value = 4626.125
$ws.Range("I73").Value = 4161.5
$ws.Range("J73").Value = 6020
$ws.Range("K73").Value = 4161.5
$ws.Range("L73").Value = 6020
$ws.Range("M73").Value = -3225.5
$ws.Range("N73").Value = -7892

$ws.Range("H97").Value = 1525.95
$ws.Range("I97").Value = 1718.4286
$ws.Range("J97").Value = 1076.8334
$ws.Range("K97").Value = 1718.4286
$ws.Range("L97").Value = 1076.8334
$ws.Range("M97").Value = -1222.4286
$ws.Range("N97").Value = -2068.8334

$ws.Range("H122").Value = 2345.3547
$ws.Range("I122").Value = 1783.0555
$ws.Range("J122").Value = 3123.923
$ws.Range("K122").Value = 5349.166499999999
$ws.Range("L122").Value = 9371.769
$ws.Range("M122").Value = -2899.166499999999
$ws.Range("N122").Value = -14271.769


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1939.2174
$ws.Range("I7").Value = 1817.7646
$ws.Range("K7").Value = 1817.7646
$ws.Range("M7").Value = -1705.7646

$ws.Range("H40").Value = 2673.2
$ws.Range("I40").Value = 2763.4546
$ws.Range("J40").Value = 2425
$ws.Range("K40").Value = 2763.4546
$ws.Range("L40").Value = 2425
$ws.Range("M40").Value = -2627.4546
$ws.Range("N40").Value = -2697

$ws.Range("H61").Value = 2873.5454
$ws.Range("I61").Value = 2513.625
$ws.Range("K61").Value = 2513.625
$ws.Range("M61").Value = -2311.625

$ws.Range("H100").Value = 6603098.5
$ws.Range("I100").Value = 7483231.5
$ws.Range("K100").Value = 7483231.5
$ws.Range("M100").Value = -7482690.5

$ws.Range("H113").Value = 2873.5454
$ws.Range("I113").Value = 2513.625
$ws.Range("K113").Value = 2513.625
$ws.Range("M113").Value = -343.625

$ws.Range("H122").Value = 3419.9033
$ws.Range("I122").Value = 4672.759
$ws.Range("J122").Value = 2318.9092
$ws.Range("K122").Value = 14018.277
$ws.Range("L122").Value = 6956.7276
$ws.Range("M122").Value = -11568.277
$ws.Range("N122").Value = -11856.7276

$ws.Range("H126").Value = 1939.2174
$ws.Range("I126").Value = 1817.7646
$ws.Range("K126").Value = 5453.293799999999
$ws.Range("M126").Value = -2983.293799999999


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1496.4445
$ws.Range("I122").Value = 1468.2858
$ws.Range("K122").Value = 4404.857400000001
$ws.Range("M122").Value = -1954.857400000001

$ws.Range("H136").Value = 1354.4166
$ws.Range("I136").Value = 1397.7097
$ws.Range("J136").Value = 1086
$ws.Range("K136").Value = 4193.1291
$ws.Range("L136").Value = 3258
$ws.Range("M136").Value = -1643.1291
$ws.Range("N136").Value = -8358

